# Apply the "search_individuals" row addition to the common_translations sheet,
# plus the small bookkeeping updates (defined name range, autofilter range,
# selection) that accompany it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("common_translations")
$ws.Activate()

# Insert a new row above the current row 66 ("search_individuals_title" /
# "Individual Advanced Search"), which pushes rows 66-77 down to 67-78.
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row with the new translation key/value pair.
$ws.Range("A66").Value = "search_individuals"
$ws.Range("B66").Value = "Search for Individuals"

# The used range grew by one row (A1:F77 -> A1:F78); re-apply AutoFilter over
# the new extent (toggle off, then back on with the new range so the filter
# descriptor in the sheet is rewritten rather than just left stale).
$ws.Range("A1:F78").AutoFilter() | Out-Null
$ws.Range("A1:F78").AutoFilter() | Out-Null

# Keep the workbook-level hidden _FilterDatabase name in sync with the new
# autofilter extent.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=common_translations!`$A`$1:`$F`$78"
    }
}

# Match the author's final selection/scroll position on this sheet.
$ws.Range("B66").Select()

"done"
